# Refresh the cryptocurrency price/volume table with the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Force the literal text into the cell (leading apostrophe = "treat
    # as text"), so numeric-looking strings like "572.49" are not auto-
    # converted into floating point numbers and reformatted by Excel.
    $range.Value = "'" + $text
}

# Row 2
$ws.Range("D2").Value = "66.501.65"
$ws.Range("E2").Value = "  -3.36%  "

# Row 3
$ws.Range("D3").Value = "3.323.50"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
Set-TextCell $ws.Range("D5") "572.49"
$ws.Range("E5").Value = "  -2.52%  "

# Row 6
Set-TextCell $ws.Range("D6") "182.98"
$ws.Range("E6").Value = "  -2.85%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.603"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "3.315.70"
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("E10").Value = "  -1.22%  "

# Row 11
Set-TextCell $ws.Range("D11") "6.64"
$ws.Range("E11").Value = "  -0.47%  "

# Row 12
Set-TextCell $ws.Range("D12") "0.404"
$ws.Range("E12").Value = "  -2.14%  "

# Row 13
$ws.Range("D13").Value = "3.896.33"
$ws.Range("E13").Value = "  +0.82%  "

# Row 14
$ws.Range("E14").Value = "  -0.91%  "

# Row 15
Set-TextCell $ws.Range("D15") "27.15"
$ws.Range("E15").Value = "  -1.81%  "

# Row 16
$ws.Range("D16").Value = "66.652.60"
$ws.Range("E16").Value = "  -3.13%  "

# Row 17
$ws.Range("E17").Value = "  -1.57%  "

# Row 18
$ws.Range("D18").Value = "3.313.83"
$ws.Range("E18").Value = "  +0.41%  "

# Row 19
Set-TextCell $ws.Range("D19") "442.87"
$ws.Range("E19").Value = "  +5.09%  "

# Row 20
Set-TextCell $ws.Range("D20") "13.57"
$ws.Range("E20").Value = "  +0.17%  "

# Row 21
Set-TextCell $ws.Range("D21") "5.67"
$ws.Range("E21").Value = "  -1.73%  "

# Row 22
Set-TextCell $ws.Range("D22") "7.67"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23
Set-TextCell $ws.Range("D23") "74.26"
$ws.Range("E23").Value = "  +3.44%  "

# Row 24
Set-TextCell $ws.Range("D24") "0.997"
$ws.Range("E24").Value = "  -0.22%  "

# Row 25
$ws.Range("D25").Value = "3.461.65"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
Set-TextCell $ws.Range("D26") "0.514"
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
Set-TextCell $ws.Range("D27") "0.0000119"
$ws.Range("E27").Value = "  +0.02%  "

# Row 28
Set-TextCell $ws.Range("D28") "0.192"
$ws.Range("E28").Value = "  +1.08%  "

# Row 29
Set-TextCell $ws.Range("D29") "8.95"
$ws.Range("E29").Value = "  -5.82%  "

# Row 30
$ws.Range("E30").Value = "  -0.89%  "

# Row 31
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
Set-TextCell $ws.Range("D32") "22.86"
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D33") "5.30"
$ws.Range("E33").Value = "  -4.07%  "

# Row 34
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws.Range("D34") "1.00"
$ws.Range("E34").Value = "  +0.09%  "

# Row 35
Set-TextCell $ws.Range("D35") "6.80"
$ws.Range("E35").Value = "  -1.83%  "

# Row 36
Set-TextCell $ws.Range("D36") "1.22"
$ws.Range("E36").Value = "  -2.91%  "

# Row 37
Set-TextCell $ws.Range("D37") "1.50"
$ws.Range("E37").Value = "  +2.07%  "

# Row 38
Set-TextCell $ws.Range("D38") "160.31"
$ws.Range("E38").Value = "  -2.45%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D39") "1.85"
$ws.Range("E39").Value = "  -3.84%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D40") "27.32"
$ws.Range("E40").Value = "  +2.51%  "

# Row 41
$ws.Range("D41").Value = "2.828.28"
$ws.Range("E41").Value = "  +5.60%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D42") "4.47"
$ws.Range("E42").Value = "  -1.33%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D43") "0.784"
$ws.Range("E43").Value = "  -1.79%  "

# Row 44
Set-TextCell $ws.Range("D44") "6.24"
$ws.Range("E44").Value = "  -2.39%  "

# Row 45
Set-TextCell $ws.Range("D45") "40.24"
$ws.Range("E45").Value = "  -0.97%  "

# Row 46
Set-TextCell $ws.Range("D46") "0.0673"
$ws.Range("E46").Value = "  -1.14%  "

# Row 47
Set-TextCell $ws.Range("D47") "24.25"
$ws.Range("E47").Value = "  -2.28%  "

# Row 48
Set-TextCell $ws.Range("D48") "2.33"
$ws.Range("E48").Value = "  -5.95%  "

# Row 49
Set-TextCell $ws.Range("D49") "321.24"
$ws.Range("E49").Value = "  -5.80%  "

# Row 50
$ws.Range("E50").Value = "  -1.97%  "

# Row 51
Set-TextCell $ws.Range("D51") "0.983"
$ws.Range("E51").Value = "  -1.56%  "
